$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(68).Insert()

$ws.Cells.Item(68,1).Value = 1
$ws.Cells.Item(68,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68,3).Value = "Arica y Parinacota"
$ws.Cells.Item(68,4).Value = 44777
$ws.Cells.Item(68,5).Value = 15
$ws.Cells.Item(68,6).Value = "Fruta"
$ws.Cells.Item(68,7).Value = 100102
$ws.Cells.Item(68,8).Value = "Cítricos"
$ws.Cells.Item(68,9).Value = 100102004
$ws.Cells.Item(68,10).Value = "Mandarina"
$ws.Cells.Item(68,11).Value = "Murcott"
$ws.Cells.Item(68,12).Value = "Segunda"
$ws.Cells.Item(68,13).Value = 250
$ws.Cells.Item(68,14).Value = 13000
$ws.Cells.Item(68,15).Value = 14000
$ws.Cells.Item(68,16).Value = 13500
$ws.Cells.Item(68,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(68,18).Value = "Región de Coquimbo"
$ws.Cells.Item(68,19).Value = 675
$ws.Cells.Item(68,20).Value = 20

Write-Host "Done"
